$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "A 11405-2024"
$ws.Range("B2").Value = 45372
$ws.Range("C2").Value = 46074
$ws.Range("F2").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G2").Value = 11.6

$ws.Range("A3").Value = "A 67756-2021"
$ws.Range("B3").Value = 44524
$ws.Range("C3").Value = 46074
$ws.Range("F3").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G3").Value = 1.6

$ws.Range("A4").Value = "A 39958-2024"
$ws.Range("B4").Value = 45553
$ws.Range("C4").Value = 46074
$ws.Range("F4").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G4").Value = 3.4

$ws.Range("A5").Value = "A 39924-2025"
$ws.Range("B5").Value = 45891
$ws.Range("C5").Value = 46074
$ws.Range("F5").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G5").Value = 1.1

$ws.Range("A6").Value = "A 39928-2025"
$ws.Range("B6").Value = 45891
$ws.Range("C6").Value = 46074
$ws.Range("F6").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G6").Value = 0.9

$ws.Range("A7").Value = "A 40001-2025"
$ws.Range("B7").Value = 45891
$ws.Range("C7").Value = 46074
$ws.Range("F7").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G7").Value = 0.6

$ws.Range("A8").Value = "A 32256-2025"
$ws.Range("B8").Value = 45835.6353125
$ws.Range("C8").Value = 46074
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = 5.4

$ws.Range("A9").Value = "A 35838-2023"
$ws.Range("B9").Value = 45147
$ws.Range("C9").Value = 46074
$ws.Range("F9").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G9").Value = 1.1

$ws.Range("A10").Value = "A 5968-2023"
$ws.Range("B10").Value = 44959
$ws.Range("C10").Value = 46074
$ws.Range("F10").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G10").Value = 1.5

$ws.Range("A11").Value = "A 2229-2023"
$ws.Range("B11").Value = 44939
$ws.Range("C11").Value = 46074
$ws.Range("F11").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G11").Value = 4.3

$ws.Range("A12").Value = "A 6004-2026"
$ws.Range("B12").Value = 46050
$ws.Range("C12").Value = 46074
$ws.Range("F12").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G12").Value = 2.7

$ws.Range("A13").Value = "A 5528-2023"
$ws.Range("B13").Value = 44957
$ws.Range("C13").Value = 46074
$ws.Range("F13").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G13").Value = 1.2

$ws.Range("A14").Value = "A 7694-2023"
$ws.Range("B14").Value = 44967
$ws.Range("C14").Value = 46074
$ws.Range("F14").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G14").Value = 2.2

$ws.Range("A15").Value = "A 2727-2024"
$ws.Range("B15").Value = 45314
$ws.Range("C15").Value = 46074
$ws.Range("F15").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G15").Value = 3.8

$ws.Range("A16").Value = "A 39876-2024"
$ws.Range("B16").Value = 45553
$ws.Range("C16").Value = 46074
$ws.Range("F16").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G16").Value = 0.3

$ws.Range("A17").Value = "A 34926-2022"
$ws.Range("B17").Value = 44796
$ws.Range("C17").Value = 46074
$ws.Range("F17").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G17").Value = 1.3

$ws.Range("A18").Value = "A 28815-2024"
$ws.Range("B18").Value = 45478
$ws.Range("C18").Value = 46074
$ws.Range("F18").Value = "Övriga statliga verk och myndigheter"
$ws.Range("G18").Value = 2.8
